# Append: 2025-10-07 18:33 JST
# Update the "取得日時" (acquired-at) timestamp for all data rows, and
# refresh row 8 / row 9 with the latest scraped data (they swapped order,
# and row 8's listing now shows the updated title/price/url that used to
# belong to row 9, and vice versa).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-07 18:33:03"

# Determine the last used row in column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Row 8 now holds what used to be row 9's listing data.
$ws.Cells.Item(8, 2).Value = "2026年度新入社員研修Javaサブ講師 (4~6月)"
$ws.Cells.Item(8, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5408522"

# Row 9 now holds what used to be row 8's listing data.
$ws.Cells.Item(9, 2).Value = "2026年度新入社員研修Javaサブ講師"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5408524"
